$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '72.180.62'
$ws.Range('E2').Value = '  +0.23%  '
$ws.Range('D3').Value = '4.044.48'
$ws.Range('E3').Value = '  -0.10%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '539.24'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '149.72'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.95%  '
$ws.Range('D7').Value = '4.035.64'
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.697'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.30%  '
$ws.Range('E9').Value = '  -0.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.752'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.65%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.172'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.26%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '53.48'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +10.54%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000334'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.61%  '
$ws.Range('E14').Value = '  -0.20%  '
$ws.Range('D15').Value = '4.685.39'
$ws.Range('D16').Value = '4.042.69'
$ws.Range('E16').Value = '  +0.17%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.33'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.40%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '20.57'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.93%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.21'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.72%  '
$ws.Range('E20').Value = '  -0.96%  '
$ws.Range('D21').Value = '72.143.45'
$ws.Range('E21').Value = '  +0.35%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '441.18'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.25%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '97.60'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.94%  '
$ws.Range('E24').Value = '  -2.30%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.26'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.11%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '14.64'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.77%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '4.48'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +23.35%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.20'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.73%  '
$ws.Range('E30').Value = '  +2.00%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '37.22'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.12%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.29'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +19.52%  '
$ws.Range('E33').Value = '  +1.65%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '13.50'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.36%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '49.35'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +13.20%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '683.63'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.10%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '67.06'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.41%  '
$ws.Range('B38').Value = 'PEPE'
$ws.Range('C38').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D38').Value = '0.0₃0915'
$ws.Range('E38').Value = '  +8.36%  '
$ws.Range('B39').Value = 'TheGraph'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.456'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.03%  '
$ws.Range('E40').Value = '  -6.03%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '11.32'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +18.23%  '
$ws.Range('B42').Value = 'ThetaToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.39'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.55%  '
$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.37'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E44').Value = '  +0.27%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.999'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.10%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0492'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.01%  '
$ws.Range('E47').Value = '  -1.14%  '
$ws.Range('E48').Value = '  -3.33%  '
$ws.Range('B49').Value = 'ApeXProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.39'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.74%  '
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.12'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.12%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.000280'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.51%  '
